$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 82919.5
$ws.Range("J3").Value = 82919.5
$ws.Range("L3").Value = 82919.5
$ws.Range("N3").Value = -83147.5
$ws.Range("H17").Value = 1147.1538
$ws.Range("J17").Value = 1147.1538
$ws.Range("L17").Value = 3441.4614
$ws.Range("N17").Value = -3777.4614
$ws.Range("H28").Value = 300.5
$ws.Range("I28").Value = 316.75
$ws.Range("J28").Value = 203
$ws.Range("K28").Value = 316.75
$ws.Range("L28").Value = 203
$ws.Range("M28").Value = 168.25
$ws.Range("N28").Value = -1173
$ws.Range("H33").Value = 1437
$ws.Range("I33").Value = 1437
$ws.Range("K33").Value = 1437
$ws.Range("M33").Value = -1208
$ws.Range("H40").Value = 11507.091
$ws.Range("J40").Value = 3158.6
$ws.Range("L40").Value = 3158.6
$ws.Range("N40").Value = -3508.6
$ws.Range("H43").Value = 685883.5
$ws.Range("I43").Value = 1501
$ws.Range("K43").Value = 1501
$ws.Range("M43").Value = -1432
$ws.Range("H64").Value = 6184
$ws.Range("I64").Value = 5651
$ws.Range("K64").Value = 5651
$ws.Range("M64").Value = -5403
$ws.Range("H67").Value = 6184
$ws.Range("I67").Value = 5651
$ws.Range("K67").Value = 5651
$ws.Range("M67").Value = -4793
$ws.Range("H80").Value = 46578.184
$ws.Range("I80").Value = 20182.2
$ws.Range("K80").Value = 60546.60000000001
$ws.Range("M80").Value = -59548.60000000001
$ws.Range("H83").Value = 46578.184
$ws.Range("I83").Value = 20182.2
$ws.Range("K83").Value = 181639.8
$ws.Range("M83").Value = -176647.8
$ws.Range("H86").Value = 75001624
$ws.Range("I86").Value = 90279380
$ws.Range("J86").Value = 13890613
$ws.Range("K86").Value = 90279380
$ws.Range("L86").Value = 13890613
$ws.Range("M86").Value = -90278257
$ws.Range("N86").Value = -13892859
$ws.Range("H89").Value = 75001624
$ws.Range("I89").Value = 90279380
$ws.Range("J89").Value = 13890613
$ws.Range("K89").Value = 451396900
$ws.Range("L89").Value = 69453065
$ws.Range("M89").Value = -451391284
$ws.Range("N89").Value = -69464297
$ws.Range("H92").Value = 3832.8333
$ws.Range("I92").Value = 699
$ws.Range("J92").Value = 6966.6665
$ws.Range("K92").Value = 699
$ws.Range("L92").Value = 6966.6665
$ws.Range("M92").Value = 549
$ws.Range("N92").Value = -9462.666499999999
$ws.Range("H98").Value = 58828976
$ws.Range("I98").Value = 66671650
$ws.Range("K98").Value = 66671650
$ws.Range("M98").Value = -66670152
$ws.Range("H99").Value = 495.75
$ws.Range("I99").Value = 594.3333
$ws.Range("K99").Value = 1782.9999
$ws.Range("M99").Value = -284.9999
$ws.Range("H101").Value = 1380.5714
$ws.Range("I101").Value = 325
$ws.Range("K101").Value = 975
$ws.Range("M101").Value = 647
$ws.Range("H102").Value = 82919.5
$ws.Range("J102").Value = 82919.5
$ws.Range("L102").Value = 82919.5
$ws.Range("N102").Value = -89409.5
$ws.Range("H106").Value = 333335100
$ws.Range("I106").Value = 333335100
$ws.Range("K106").Value = 333335100
$ws.Range("M106").Value = -333334469
$ws.Range("H107").Value = 62503570
$ws.Range("I107").Value = 41669136
$ws.Range("K107").Value = 41669136
$ws.Range("M107").Value = -41667216
$ws.Range("H116").Value = 9620014
$ws.Range("I116").Value = 25002400
$ws.Range("J116").Value = 6022.875
$ws.Range("K116").Value = 25002400
$ws.Range("L116").Value = 6022.875
$ws.Range("M116").Value = -24998958
$ws.Range("N116").Value = -12906.875
$ws.Range("H122").Value = 58828976
$ws.Range("I122").Value = 66671650
$ws.Range("K122").Value = 200014950
$ws.Range("M122").Value = -200012500
$ws.Range("H138").Value = 5754.1567
$ws.Range("I138").Value = 2485.6667
$ws.Range("J138").Value = 7116.028
$ws.Range("K138").Value = 7457.000100000001
$ws.Range("L138").Value = 21348.084
$ws.Range("M138").Value = -2317.000100000001
$ws.Range("N138").Value = -31628.084

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2122176.8
$ws.Range("I32").Value = 2361970.5
$ws.Range("J32").Value = 3999
$ws.Range("K32").Value = 2361970.5
$ws.Range("L32").Value = 3999
$ws.Range("M32").Value = -2361683.5
$ws.Range("N32").Value = -4573
$ws.Range("H61").Value = 4391.946
$ws.Range("I61").Value = 1287.9395
$ws.Range("K61").Value = 1287.9395
$ws.Range("M61").Value = -1075.9395
$ws.Range("H102").Value = 4639.4
$ws.Range("I102").Value = 4499.25
$ws.Range("K102").Value = 4499.25
$ws.Range("M102").Value = -2877.25
$ws.Range("H122").Value = 11037.869
$ws.Range("I122").Value = 14200.5
$ws.Range("K122").Value = 42601.5
$ws.Range("M122").Value = -40151.5
$ws.Range("H132").Value = 4690.68
$ws.Range("I132").Value = 2042.0312
$ws.Range("J132").Value = 9399.388999999999
$ws.Range("K132").Value = 6126.0936
$ws.Range("L132").Value = 28198.167
$ws.Range("M132").Value = -3596.0936
$ws.Range("N132").Value = -33258.167
$ws.Range("H136").Value = 4391.946
$ws.Range("I136").Value = 1287.9395
$ws.Range("K136").Value = 3863.8185
$ws.Range("M136").Value = -1313.8185

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10418126
$ws.Range("I20").Value = 16668354
$ws.Range("J20").Value = 1080
$ws.Range("K20").Value = 16668354
$ws.Range("L20").Value = 1080
$ws.Range("M20").Value = -16668107
$ws.Range("N20").Value = -1574
$ws.Range("H22").Value = 16914.334
$ws.Range("I22").Value = 16914.334
$ws.Range("K22").Value = 16914.334
$ws.Range("M22").Value = -16741.334
$ws.Range("H28").Value = 48271
$ws.Range("J28").Value = 48271
$ws.Range("L28").Value = 48271
$ws.Range("N28").Value = -48859
$ws.Range("H40").Value = 45397.5
$ws.Range("J40").Value = 45397.5
$ws.Range("L40").Value = 45397.5
$ws.Range("N40").Value = -45927.5
$ws.Range("H99").Value = 8265919.5
$ws.Range("I99").Value = 883.1429000000001
$ws.Range("K99").Value = 883.1429000000001
$ws.Range("M99").Value = 614.8570999999999
$ws.Range("H105").Value = 2978.4783
$ws.Range("I105").Value = 2142.8572
$ws.Range("K105").Value = 2142.8572
$ws.Range("M105").Value = -395.8571999999999
$ws.Range("H134").Value = 6807.1177
$ws.Range("I134").Value = 3302.842
$ws.Range("K134").Value = 9908.526
$ws.Range("M134").Value = -7373.526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 45455188
$ws.Range("I7").Value = 1035
$ws.Range("J7").Value = 83333650
$ws.Range("K7").Value = 1035
$ws.Range("L7").Value = 83333650
$ws.Range("M7").Value = -922
$ws.Range("N7").Value = -83333876
$ws.Range("H31").Value = 9118.052
$ws.Range("I31").Value = 4017.4736
$ws.Range("J31").Value = 13963.6
$ws.Range("K31").Value = 4017.4736
$ws.Range("L31").Value = 13963.6
$ws.Range("M31").Value = -3722.4736
$ws.Range("N31").Value = -14553.6
$ws.Range("H34").Value = 9118.052
$ws.Range("I34").Value = 4017.4736
$ws.Range("J34").Value = 13963.6
$ws.Range("K34").Value = 4017.4736
$ws.Range("L34").Value = 13963.6
$ws.Range("M34").Value = -3815.4736
$ws.Range("N34").Value = -14367.6
$ws.Range("H43").Value = 37578
$ws.Range("J43").Value = 37578
$ws.Range("L43").Value = 37578
$ws.Range("N43").Value = -37946
$ws.Range("H86").Value = 56835220
$ws.Range("I86").Value = 24267702
$ws.Range("J86").Value = 111114420
$ws.Range("K86").Value = 24267702
$ws.Range("L86").Value = 111114420
$ws.Range("M86").Value = -24266579
$ws.Range("N86").Value = -111116666
$ws.Range("H89").Value = 56835220
$ws.Range("I89").Value = 24267702
$ws.Range("J89").Value = 111114420
$ws.Range("K89").Value = 121338510
$ws.Range("L89").Value = 555572100
$ws.Range("M89").Value = -121332894
$ws.Range("N89").Value = -555583332
$ws.Range("H99").Value = 9473.556
$ws.Range("I99").Value = 12114.143
$ws.Range("K99").Value = 12114.143
$ws.Range("M99").Value = -10616.143
$ws.Range("H101").Value = 37578
$ws.Range("J101").Value = 37578
$ws.Range("L101").Value = 37578
$ws.Range("N101").Value = -44068
$ws.Range("H120").Value = 80000
$ws.Range("J120").Value = 80000
$ws.Range("L120").Value = 80000
$ws.Range("N120").Value = -87258
$ws.Range("H122").Value = 2998
$ws.Range("I122").Value = 2998
$ws.Range("K122").Value = 8994
$ws.Range("M122").Value = -6544
$ws.Range("H125").Value = 51250.5
$ws.Range("J125").Value = 51250.5
$ws.Range("L125").Value = 51250.5
$ws.Range("N125").Value = -56170.5
$ws.Range("H126").Value = 9473.556
$ws.Range("I126").Value = 12114.143
$ws.Range("K126").Value = 36342.429
$ws.Range("M126").Value = -33872.429
$ws.Range("H132").Value = 6579.2104
$ws.Range("I132").Value = 2317.889
$ws.Range("J132").Value = 10414.4
$ws.Range("K132").Value = 6953.667
$ws.Range("L132").Value = 31243.2
$ws.Range("M132").Value = -4423.667
$ws.Range("N132").Value = -36303.2
$ws.Range("H134").Value = 7839.711
$ws.Range("I134").Value = 8008.048
$ws.Range("J134").Value = 7692.4165
$ws.Range("K134").Value = 24024.144
$ws.Range("L134").Value = 23077.2495
$ws.Range("M134").Value = -21489.144
$ws.Range("N134").Value = -28147.2495
$ws.Range("H141").Value = 70999.5
$ws.Range("J141").Value = 70999.5
$ws.Range("L141").Value = 70999.5
$ws.Range("N141").Value = -81359.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1729.9
$ws.Range("I7").Value = 274.75
$ws.Range("J7").Value = 2700
$ws.Range("K7").Value = 824.25
$ws.Range("L7").Value = 8100
$ws.Range("M7").Value = -712.25
$ws.Range("N7").Value = -8324
$ws.Range("H12").Value = 2381724.2
$ws.Range("J12").Value = 3846795.5
$ws.Range("L12").Value = 11540386.5
$ws.Range("N12").Value = -11540732.5
$ws.Range("H98").Value = 83333580
$ws.Range("J98").Value = 166666830
$ws.Range("L98").Value = 500000490
$ws.Range("N98").Value = -500003486
$ws.Range("H113").Value = 2838.5715
$ws.Range("I113").Value = 1388.75
$ws.Range("J113").Value = 3179.7058
$ws.Range("K113").Value = 4166.25
$ws.Range("L113").Value = 9539.117400000001
$ws.Range("M113").Value = -1996.25
$ws.Range("N113").Value = -13879.1174
$ws.Range("H137").Value = 251495.38
$ws.Range("I137").Value = 167827.33
$ws.Range("K137").Value = 503481.99
$ws.Range("M137").Value = -498381.99

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2815.0908
$ws.Range("I80").Value = 2995.5557
$ws.Range("K80").Value = 2995.5557
$ws.Range("M80").Value = -1997.5557
$ws.Range("H83").Value = 2815.0908
$ws.Range("I83").Value = 2995.5557
$ws.Range("K83").Value = 14977.7785
$ws.Range("M83").Value = -9985.7785
$ws.Range("H122").Value = 50519.617
$ws.Range("I122").Value = 68434.47
$ws.Range("K122").Value = 205303.41
$ws.Range("M122").Value = -202853.41

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7000
$ws.Range("J7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("N7").Value = -7224
$ws.Range("H16").Value = 1058.7858
$ws.Range("I16").Value = 1151.0834
$ws.Range("K16").Value = 1151.0834
$ws.Range("M16").Value = -981.0834
$ws.Range("H22").Value = 29856.428
$ws.Range("I22").Value = 7998.75
$ws.Range("J22").Value = 59000
$ws.Range("K22").Value = 7998.75
$ws.Range("L22").Value = 59000
$ws.Range("M22").Value = -7703.75
$ws.Range("N22").Value = -59590
$ws.Range("H27").Value = 29856.428
$ws.Range("I27").Value = 7998.75
$ws.Range("J27").Value = 59000
$ws.Range("K27").Value = 7998.75
$ws.Range("L27").Value = 59000
$ws.Range("M27").Value = -7891.75
$ws.Range("N27").Value = -59214
$ws.Range("H40").Value = 33337386
$ws.Range("I40").Value = 55559150
$ws.Range("J40").Value = 4734.1665
$ws.Range("K40").Value = 55559150
$ws.Range("L40").Value = 4734.1665
$ws.Range("M40").Value = -55559014
$ws.Range("N40").Value = -5006.1665
$ws.Range("H46").Value = 12346676
$ws.Range("I46").Value = 895
$ws.Range("K46").Value = 895
$ws.Range("M46").Value = -707
$ws.Range("H59").Value = 56546
$ws.Range("J59").Value = 56546
$ws.Range("L59").Value = 56546
$ws.Range("N59").Value = -57854
$ws.Range("H61").Value = 5439.0625
$ws.Range("I61").Value = 4098.8
$ws.Range("J61").Value = 6048.273
$ws.Range("K61").Value = 4098.8
$ws.Range("L61").Value = 6048.273
$ws.Range("M61").Value = -3896.8
$ws.Range("N61").Value = -6452.273
$ws.Range("H113").Value = 5439.0625
$ws.Range("I113").Value = 4098.8
$ws.Range("J113").Value = 6048.273
$ws.Range("K113").Value = 4098.8
$ws.Range("L113").Value = 6048.273
$ws.Range("M113").Value = -1928.8
$ws.Range("N113").Value = -10388.273
$ws.Range("H122").Value = 3519.9111
$ws.Range("I122").Value = 2797
$ws.Range("J122").Value = 4830.1875
$ws.Range("K122").Value = 8391
$ws.Range("L122").Value = 14490.5625
$ws.Range("M122").Value = -5941
$ws.Range("N122").Value = -19390.5625
$ws.Range("H126").Value = 7000
$ws.Range("J126").Value = 7000
$ws.Range("L126").Value = 21000
$ws.Range("N126").Value = -25940
$ws.Range("H132").Value = 16674700
$ws.Range("I132").Value = 31254562
$ws.Range("J132").Value = 12000.286
$ws.Range("K132").Value = 93763686
$ws.Range("L132").Value = 36000.858
$ws.Range("M132").Value = -93761156
$ws.Range("N132").Value = -41060.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 38888
$ws.Range("J51").Value = 38888
$ws.Range("L51").Value = 38888
$ws.Range("N51").Value = -39908
$ws.Range("H101").Value = 30578.8
$ws.Range("J101").Value = 30578.8
$ws.Range("L101").Value = 30578.8
$ws.Range("N101").Value = -37068.8
$ws.Range("H113").Value = 1735.7693
$ws.Range("I113").Value = 1537.8
$ws.Range("J113").Value = 2005.7273
$ws.Range("K113").Value = 4613.4
$ws.Range("L113").Value = 6017.1819
$ws.Range("M113").Value = -2443.4
$ws.Range("N113").Value = -10357.1819
$ws.Range("H122").Value = 19388474
$ws.Range("I122").Value = 33602970
$ws.Range("J122").Value = 5076.636
$ws.Range("K122").Value = 100808910
$ws.Range("L122").Value = 15229.908
$ws.Range("M122").Value = -100806460
$ws.Range("N122").Value = -20129.908
$ws.Range("H132").Value = 62576624
$ws.Range("I132").Value = 100022200
$ws.Range("J132").Value = 167333.33
$ws.Range("K132").Value = 300066600
$ws.Range("M132").Value = -300064070
